{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst p = paragraphs.items[0];\n\n// \"Versi\" + \"on\" (two runs) -> single run \"Version\" (merges the split word\n// back together, matching the diff's collapsed run).\nconst versionHits = p.search(\"Version\", { matchCase: true });\nversionHits.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < versionHits.items.length; i++) {\n  versionHits.items[i].insertText(\"Version\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// \" 2\" + \".\" (the \".\" trails a bookmark in its own run) -> \" 1.\" merged\n// into a single run, dropping the now-empty trailing run.\nconst numberHits = p.search(\"2.\", { matchCase: true });\nnumberHits.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < numberHits.items.length; i++) {\n  numberHits.items[i].insertText(\"1.\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# \"Versi\" + \"on\" (two runs) -> single run \"Version\" (merges the split word\n# back together, matching the diff's collapsed run).\n$find = $d.Content.Find\n$find.Text = \"Version\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, \"Version\", 2)\n\n# Locate \"2.\" without replacing yet. The \"2\" and the \".\" are in separate\n# runs straddling the \"_GoBack\" bookmark, so we capture the offsets first.\n$hit = $d.Range(0, 0)\n$findDot = $hit.Find\n$findDot.Text = \"2.\"\n$findDot.Execute($findDot.Text, $false, $false, $false, $false, $false, $true, 0, $false, \"\", 0)\n$start = $hit.Start\n$end = $hit.End\n\n# The trailing \".\" sits in its own run just after the bookmark. Remove it\n# first so the next edit doesn't have to span the bookmark (which would\n# otherwise delete the bookmark).\n$dot = $d.Range($end - 1, $end)\n$dot.Delete()\n\n# Replace the \"2\" (now immediately before the bookmark) with \"1.\" so it\n# merges into the \" 2\" run, producing \" 1.\" and leaving the bookmark intact.\n$two = $d.Range($start, $start + 1)\n$two.Text = \"1.\"\n"}
